$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ligand/Receptor-expressing cell counts (1 -> 3) and all downstream
# specificity/expression statistics recomputed from the new cell counts,
# per Dr Hou's advice.
$updates = @{
    2 = @{ 'E' = 3; 'G' = 203.7816646666667; 'H' = 611.344994; 'I' = 0.6667327591988204; 'J' = 0.6667327591988205; 'K' = 3; 'M' = 23.18731733333334; 'N' = 69.561952; 'O' = 0.8148739324305957; 'P' = 0.8148739324305956; 'Q' = 4725.150125340921; 'R' = 42526.35112806829; 'S' = 0.5433031453686442; 'T' = 0.5433031453686442 }
    3 = @{ 'E' = 3; 'G' = 203.7816646666667; 'H' = 611.344994; 'I' = 0.6667327591988204; 'J' = 0.6667327591988205; 'K' = 3; 'M' = 5.267780333333334; 'N' = 15.803341; 'O' = 0.1851260675694043; 'P' = 0.1851260675694043; 'Q' = 1073.477045424995; 'R' = 9661.293408824955; 'S' = 0.1234296138301762; 'T' = 0.1234296138301762 }
    4 = @{ 'E' = 3; 'G' = 63.14058933333333; 'H' = 189.421768; 'I' = 0.2065833519051582; 'J' = 0.2065833519051582; 'K' = 3; 'M' = 23.18731733333334; 'N' = 69.561952; 'O' = 0.8148739324305957; 'P' = 0.8148739324305956; 'Q' = 1464.060881485682; 'R' = 13176.54793337114; 'S' = 0.1683393883416498; 'T' = 0.1683393883416498 }
    5 = @{ 'E' = 3; 'G' = 63.14058933333333; 'H' = 189.421768; 'I' = 0.2065833519051582; 'J' = 0.2065833519051582; 'K' = 3; 'M' = 5.267780333333334; 'N' = 15.803341; 'O' = 0.1851260675694043; 'P' = 0.1851260675694043; 'Q' = 332.6107547252098; 'R' = 2993.496792526888; 'S' = 0.03824396356350834; 'T' = 0.03824396356350835 }
    6 = @{ 'E' = 3; 'G' = 38.719942; 'H' = 116.159826; 'I' = 0.1266838888960214; 'J' = 0.1266838888960214; 'K' = 3; 'M' = 23.18731733333334; 'N' = 69.561952; 'O' = 0.8148739324305957; 'P' = 0.8148739324305956; 'Q' = 897.8115822822614; 'R' = 8080.304240540352; 'S' = 0.1032313987203016; 'T' = 0.1032313987203016 }
    7 = @{ 'E' = 3; 'G' = 38.719942; 'H' = 116.159826; 'I' = 0.1266838888960214; 'J' = 0.1266838888960214; 'K' = 3; 'M' = 5.267780333333334; 'N' = 15.803341; 'O' = 0.1851260675694043; 'P' = 0.1851260675694043; 'Q' = 203.9681489754073; 'R' = 1835.713340778666; 'S' = 0.02345249017571977; 'T' = 0.02345249017571977 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
